$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2703.7144
$ws.Range("I70").Value = 2703.7144
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 8111.1432
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -7841.1432
$ws.Range("N70").Value = ""
$ws.Range("H73").Value = 2703.7144
$ws.Range("I73").Value = 2703.7144
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 8111.1432
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -7175.1432
$ws.Range("N73").Value = ""
$ws.Range("H97").Value = 2496.55
$ws.Range("J97").Value = 2496.55
$ws.Range("L97").Value = 7489.650000000001
$ws.Range("N97").Value = -8481.650000000001
$ws.Range("H98").Value = 3458.8572
$ws.Range("I98").Value = 3618.6667
$ws.Range("J98").Value = 2500
$ws.Range("K98").Value = 3618.6667
$ws.Range("L98").Value = 2500
$ws.Range("M98").Value = -2120.6667
$ws.Range("N98").Value = -5496
$ws.Range("H122").Value = 3458.8572
$ws.Range("I122").Value = 3618.6667
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 10856.0001
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -8406.000100000001
$ws.Range("N122").Value = -12400
$ws.Range("H132").Value = 6927.1
$ws.Range("I132").Value = 3681.5715
$ws.Range("K132").Value = 11044.7145
$ws.Range("M132").Value = -8514.7145
$ws.Range("H137").Value = 1136.2059
$ws.Range("I137").Value = 974.1579
$ws.Range("K137").Value = 2922.4737
$ws.Range("M137").Value = -372.4737
$ws.Range("H138").Value = 1931.2094
$ws.Range("I138").Value = 1335
$ws.Range("J138").Value = 2089.0293
$ws.Range("K138").Value = 4005
$ws.Range("L138").Value = 6267.0879
$ws.Range("M138").Value = 1135
$ws.Range("N138").Value = -16547.0879

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 922.7
$ws.Range("I61").Value = 691.8889
$ws.Range("K61").Value = 691.8889
$ws.Range("M61").Value = -479.8889
$ws.Range("H74").Value = 739.9231
$ws.Range("I74").Value = 726.13635
$ws.Range("J74").Value = 815.75
$ws.Range("K74").Value = 726.13635
$ws.Range("L74").Value = 815.75
$ws.Range("M74").Value = 147.86365
$ws.Range("N74").Value = -2563.75
$ws.Range("H77").Value = 739.9231
$ws.Range("I77").Value = 726.13635
$ws.Range("J77").Value = 815.75
$ws.Range("K77").Value = 3630.68175
$ws.Range("L77").Value = 4078.75
$ws.Range("M77").Value = 737.3182500000003
$ws.Range("N77").Value = -12814.75
$ws.Range("H132").Value = 2463.3794
$ws.Range("I132").Value = 2549.5789
$ws.Range("K132").Value = 7648.736699999999
$ws.Range("M132").Value = -5118.736699999999
$ws.Range("H136").Value = 922.7
$ws.Range("I136").Value = 691.8889
$ws.Range("K136").Value = 2075.6667
$ws.Range("M136").Value = 474.3332999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 19647.334
$ws.Range("I26").Value = 19647.334
$ws.Range("K26").Value = 19647.334
$ws.Range("M26").Value = -19355.334
$ws.Range("H94").Value = 41667668
$ws.Range("I94").Value = 62500524
$ws.Range("J94").Value = 1950
$ws.Range("K94").Value = 62500524
$ws.Range("L94").Value = 1950
$ws.Range("M94").Value = -62500073
$ws.Range("N94").Value = -2852
$ws.Range("H109").Value = 20000
$ws.Range("I109").Value = 20000
$ws.Range("K109").Value = 20000
$ws.Range("M109").Value = -18613
$ws.Range("H134").Value = 10376
$ws.Range("I134").Value = 1939.25
$ws.Range("J134").Value = 27249.5
$ws.Range("K134").Value = 5817.75
$ws.Range("L134").Value = 81748.5
$ws.Range("M134").Value = -3282.75
$ws.Range("N134").Value = -86818.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1193.7567
$ws.Range("I31").Value = 726.86365
$ws.Range("J31").Value = 1878.5333
$ws.Range("K31").Value = 726.86365
$ws.Range("L31").Value = 1878.5333
$ws.Range("M31").Value = -431.86365
$ws.Range("N31").Value = -2468.5333
$ws.Range("H34").Value = 1193.7567
$ws.Range("I34").Value = 726.86365
$ws.Range("J34").Value = 1878.5333
$ws.Range("K34").Value = 726.86365
$ws.Range("L34").Value = 1878.5333
$ws.Range("M34").Value = -524.86365
$ws.Range("N34").Value = -2282.5333
$ws.Range("H41").Value = 6851.143
$ws.Range("J41").Value = 27000
$ws.Range("L41").Value = 27000
$ws.Range("N41").Value = -27856
$ws.Range("H58").Value = 801.3333
$ws.Range("I58").Value = 855.2857
$ws.Range("J58").Value = 675.44446
$ws.Range("K58").Value = 855.2857
$ws.Range("L58").Value = 675.44446
$ws.Range("M58").Value = -652.2857
$ws.Range("N58").Value = -1081.44446
$ws.Range("H99").Value = 2127
$ws.Range("I99").Value = 1987
$ws.Range("J99").Value = 2407
$ws.Range("K99").Value = 1987
$ws.Range("L99").Value = 2407
$ws.Range("M99").Value = -489
$ws.Range("N99").Value = -5403
$ws.Range("H126").Value = 2127
$ws.Range("I126").Value = 1987
$ws.Range("J126").Value = 2407
$ws.Range("K126").Value = 5961
$ws.Range("L126").Value = 7221
$ws.Range("M126").Value = -3491
$ws.Range("N126").Value = -12161
$ws.Range("H132").Value = 7116.1904
$ws.Range("I132").Value = 9232.923000000001
$ws.Range("K132").Value = 27698.769
$ws.Range("M132").Value = -25168.769
$ws.Range("H134").Value = 11112600
$ws.Range("I134").Value = 13890213
$ws.Range("K134").Value = 41670639
$ws.Range("M134").Value = -41668104
$ws.Range("H136").Value = 801.3333
$ws.Range("I136").Value = 855.2857
$ws.Range("J136").Value = 675.44446
$ws.Range("K136").Value = 2565.8571
$ws.Range("L136").Value = 2026.33338
$ws.Range("M136").Value = -15.85710000000017
$ws.Range("N136").Value = -7126.33338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1581.7778
$ws.Range("I5").Value = 1366.8572
$ws.Range("K5").Value = 4100.571599999999
$ws.Range("M5").Value = -3988.571599999999
$ws.Range("H68").Value = 1331.7894
$ws.Range("J68").Value = 1540.2142
$ws.Range("L68").Value = 4620.642599999999
$ws.Range("N68").Value = -6242.642599999999
$ws.Range("H70").Value = 4468.75
$ws.Range("J70").Value = 5383.3335
$ws.Range("L70").Value = 16150.0005
$ws.Range("N70").Value = -16780.0005
$ws.Range("H71").Value = 1331.7894
$ws.Range("J71").Value = 1540.2142
$ws.Range("L71").Value = 13861.9278
$ws.Range("N71").Value = -21973.9278
$ws.Range("H73").Value = 4468.75
$ws.Range("J73").Value = 5383.3335
$ws.Range("L73").Value = 16150.0005
$ws.Range("N73").Value = -18334.0005
$ws.Range("H92").Value = 1150
$ws.Range("J92").Value = 1150
$ws.Range("L92").Value = 3450
$ws.Range("N92").Value = -5946
$ws.Range("H98").Value = 762.8889
$ws.Range("J98").Value = 1996.6666
$ws.Range("L98").Value = 5989.9998
$ws.Range("N98").Value = -8985.9998
$ws.Range("H122").Value = 867.9167
$ws.Range("I122").Value = 498.6
$ws.Range("J122").Value = 1131.7142
$ws.Range("K122").Value = 4487.400000000001
$ws.Range("L122").Value = 10185.4278
$ws.Range("M122").Value = -2037.400000000001
$ws.Range("N122").Value = -15085.4278
$ws.Range("H135").Value = 1581.7778
$ws.Range("I135").Value = 1366.8572
$ws.Range("K135").Value = 12301.7148
$ws.Range("M135").Value = -9766.7148
$ws.Range("H137").Value = 5140.4346
$ws.Range("J137").Value = 5789
$ws.Range("L137").Value = 17367
$ws.Range("N137").Value = -27567

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 56255164
$ws.Range("I70").Value = 125002650
$ws.Range("J70").Value = 33339334
$ws.Range("K70").Value = 125002650
$ws.Range("L70").Value = 33339334
$ws.Range("M70").Value = -125002380
$ws.Range("N70").Value = -33339874
$ws.Range("H73").Value = 56255164
$ws.Range("I73").Value = 125002650
$ws.Range("J73").Value = 33339334
$ws.Range("K73").Value = 125002650
$ws.Range("L73").Value = 33339334
$ws.Range("M73").Value = -125001714
$ws.Range("N73").Value = -33341206
$ws.Range("H122").Value = 3466.5
$ws.Range("I122").Value = 2816.1667
$ws.Range("K122").Value = 8448.500100000001
$ws.Range("M122").Value = -5998.500100000001
$ws.Range("H126").Value = 2181.6667
$ws.Range("I126").Value = 1772.5
$ws.Range("K126").Value = 5317.5
$ws.Range("M126").Value = -2847.5
$ws.Range("H132").Value = 2222.6765
$ws.Range("I132").Value = 1986.55
$ws.Range("K132").Value = 5959.65
$ws.Range("M132").Value = -3429.65

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2152.5
$ws.Range("I7").Value = 2500
$ws.Range("J7").Value = 1805
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 1805
$ws.Range("M7").Value = -2388
$ws.Range("N7").Value = -2029
$ws.Range("H109").Value = 25000
$ws.Range("J109").Value = 25000
$ws.Range("L109").Value = 25000
$ws.Range("N109").Value = -27774
$ws.Range("H126").Value = 2152.5
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 1805
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 5415
$ws.Range("M126").Value = -5030
$ws.Range("N126").Value = -10355
$ws.Range("H136").Value = 2201.1538
$ws.Range("I136").Value = 2261
$ws.Range("J136").Value = 2001.6666
$ws.Range("K136").Value = 6783
$ws.Range("L136").Value = 6004.9998
$ws.Range("M136").Value = -4233
$ws.Range("N136").Value = -11104.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2785.4
$ws.Range("I132").Value = 2500.3845
$ws.Range("J132").Value = 3314.7144
$ws.Range("K132").Value = 7501.1535
$ws.Range("L132").Value = 9944.143199999999
$ws.Range("M132").Value = -4971.1535
$ws.Range("N132").Value = -15004.1432
$ws.Range("H136").Value = 691.6429000000001
$ws.Range("I136").Value = 514.1539
$ws.Range("K136").Value = 1542.4617
$ws.Range("M136").Value = 1007.5383
